$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the previous year row (A4) onto the new year-label cell (A5),
# then set its value -- mirrors the bold/centered/bordered style used for every year label.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "2021年"

# Populate the 2021 data row (row 5) across the same columns (B:DK) used by the
# prior year rows. Columns not listed here have no reported value for 2021.
$ws.Range("B5").Value = 122.9
$ws.Range("C5").Value = 25.1
$ws.Range("E5").Value = 7.2
$ws.Range("F5").Value = 23.5
$ws.Range("H5").Value = -19
$ws.Range("I5").Value = -2.1
$ws.Range("J5").Value = 1.1
$ws.Range("K5").Value = -41.7
$ws.Range("M5").Value = -6.3
$ws.Range("O5").Value = -16.8
$ws.Range("P5").Value = -8.1
$ws.Range("Q5").Value = -46.3
$ws.Range("T5").Value = -50.7
$ws.Range("U5").Value = -20.5
$ws.Range("V5").Value = -42.4
$ws.Range("W5").Value = 24.8
$ws.Range("X5").Value = 1.7
$ws.Range("Y5").Value = 17.1
$ws.Range("Z5").Value = 34.6
$ws.Range("AA5").Value = 13.5
$ws.Range("AB5").Value = -13.4
$ws.Range("AC5").Value = -27
$ws.Range("AD5").Value = 19.9
$ws.Range("AE5").Value = 43.9
$ws.Range("AK5").Value = 24.5
$ws.Range("AL5").Value = 13.1
$ws.Range("AM5").Value = 3.1
$ws.Range("AN5").Value = -84.7
$ws.Range("AO5").Value = -92.2
$ws.Range("AP5").Value = -63.9
$ws.Range("AQ5").Value = -8.3
$ws.Range("AU5").Value = 151.3
$ws.Range("AV5").Value = 13.5
$ws.Range("AX5").Value = -4.3
$ws.Range("AY5").Value = 42.3
$ws.Range("AZ5").Value = 7.8
$ws.Range("BA5").Value = -13.5
$ws.Range("BB5").Value = -50.1
$ws.Range("BC5").Value = 46
$ws.Range("BE5").Value = 8.8
$ws.Range("BF5").Value = -53.9
$ws.Range("BG5").Value = 49.3
$ws.Range("BJ5").Value = 2
$ws.Range("BK5").Value = 145.4
$ws.Range("BL5").Value = -17.7
$ws.Range("BM5").Value = -78.3
$ws.Range("BN5").Value = 11.6
$ws.Range("BO5").Value = -24.4
$ws.Range("BP5").Value = 166.9
$ws.Range("BR5").Value = 83.6
$ws.Range("BS5").Value = 2.7
$ws.Range("BT5").Value = -16.1
$ws.Range("BU5").Value = -9
$ws.Range("BV5").Value = -11.7
$ws.Range("BW5").Value = -19.9
$ws.Range("BX5").Value = 17.1
$ws.Range("BY5").Value = -16
$ws.Range("BZ5").Value = 2.5
$ws.Range("CA5").Value = -61.5
$ws.Range("CB5").Value = 288.9
$ws.Range("CC5").Value = 66.1
$ws.Range("CE5").Value = 53.3
$ws.Range("CF5").Value = 82.6
$ws.Range("CG5").Value = 127.1
$ws.Range("CH5").Value = 139.9
$ws.Range("CI5").Value = 47.4
$ws.Range("CK5").Value = 34.4
$ws.Range("CL5").Value = -0.1
$ws.Range("CO5").Value = -15.4
$ws.Range("CP5").Value = 6.9
$ws.Range("CQ5").Value = 29.6
$ws.Range("CS5").Value = -3.8
$ws.Range("CT5").Value = -1.1
$ws.Range("CU5").Value = -0.9
$ws.Range("CV5").Value = 113.2
$ws.Range("CW5").Value = -28.7
$ws.Range("CX5").Value = 3.7
$ws.Range("CY5").Value = 147.6
$ws.Range("CZ5").Value = 198
$ws.Range("DA5").Value = 9.9
$ws.Range("DB5").Value = 67
$ws.Range("DC5").Value = 3.4
$ws.Range("DD5").Value = 60.7
$ws.Range("DE5").Value = 63.6
$ws.Range("DF5").Value = 21.8
$ws.Range("DG5").Value = 287.8
$ws.Range("DH5").Value = 15.8
$ws.Range("DI5").Value = 92.7
$ws.Range("DJ5").Value = -12.1
$ws.Range("DK5").Value = -14.5
